$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new observation record ("Tvåtandad spolsnäcka") was inserted into the
# sightings table, pushing the former rows 25-28 down to rows 26-29, and the
# individual records were re-sequenced/refreshed (species re-identified,
# coordinates rounded, timestamps corrected). Insert a blank row at 25 first
# so everything below shifts down, matching the append of row 29.
# ---------------------------------------------------------------------------

$ws.Rows(25).Insert()

# The freshly inserted row has no cells at all yet; stamp the handful of
# always-blank columns (matching the shape every other record row uses) onto
# it by copying from an already-blank neighbour cell before filling in data.
$ws.Range("K26").Copy($ws.Range("I25"))
$ws.Range("K26").Copy($ws.Range("K25"))
$ws.Range("AT26").Copy($ws.Range("AT25"))
$ws.Range("AY26").Copy($ws.Range("AY25"))

# --- Row 24: "Vanlig groda" record - coordinates rounded ------------------
$ws.Range("Q24").Value = 580550
$ws.Range("R24").Value = 6579522

# --- Row 25 (brand-new row): "Mindre märgborre" record ---------------------
$ws.Range("A25").Value = 112195278
$ws.Range("B25").Value = 8377
$ws.Range("C25").Value = "Ovaliderad"
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 106545
$ws.Range("F25").Value = "Mindre märgborre"
$ws.Range("G25").Value = "Tomicus minor"
$ws.Range("H25").Value = "(Hartig, 1834)"
$ws.Range("P25").Value = "Flugmötesskogen , Srm"
$ws.Range("Q25").Value = 580550
$ws.Range("R25").Value = 6579320
$ws.Range("S25").Value = 5
$ws.Range("T25").Value = "Södermanland"
$ws.Range("U25").Value = "Eskilstuna"
$ws.Range("V25").Value = "Södermanland"
$ws.Range("W25").Value = "Eskilstuna"
# Date-like text must stay text (not become a date serial) - pull it in via
# copy/paste from an identical existing text cell instead of a raw assignment.
$ws.Range("Y26").Copy() | Out-Null
$ws.Range("Y25").PasteSpecial() | Out-Null
$ws.Range("Z25").Value = "11:00"
$ws.Range("AA26").Copy() | Out-Null
$ws.Range("AA25").PasteSpecial() | Out-Null
$ws.Range("AB25").Value = "11:00"
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AW25").Value = "Ella Axelsson Elfving"
$ws.Range("AX25").Value = "Ella Axelsson Elfving"

# --- Row 26 (was row 25): now the "Spillkråka" record -----------------------
$ws.Range("A26").Value = 112196861
$ws.Range("B26").Value = 56414
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 100049
$ws.Range("F26").Value = "Spillkråka"
$ws.Range("G26").Value = "Dryocopus martius"
$ws.Range("H26").Value = "(Linnaeus, 1758)"
$ws.Range("Q26").Value = 580550
$ws.Range("R26").Value = 6579320
$ws.Range("Z26").Value = "12:30"
$ws.Range("AB26").Value = "12:30"

# --- Row 27 (was row 26): now the "Hasselticka" record ----------------------
$ws.Range("A27").Value = 112196324
$ws.Range("B27").Value = 89953
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 3884
$ws.Range("F27").Value = "Hasselticka"
$ws.Range("G27").Value = "Dichomitus campestris"
$ws.Range("H27").Value = "(Quél.) Domański & Orlicz"
$ws.Range("Q27").Value = 580550
$ws.Range("R27").Value = 6579320
$ws.Range("Z27").Value = "12:13"
$ws.Range("AB27").Value = "12:13"

# --- Row 28 (was row 27): now the "Jättesvampmal" record; drop its comment --
$ws.Range("A28").Value = 112196967
$ws.Range("B28").Value = 43467
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 101735
$ws.Range("F28").Value = "Jättesvampmal"
$ws.Range("G28").Value = "Scardia boletella"
$ws.Range("H28").Value = "(Fabricius, 1794)"
$ws.Range("Q28").Value = 580550
$ws.Range("R28").Value = 6579320
$ws.Range("Z28").Value = "12:30"
$ws.Range("AB28").Value = "12:30"
$ws.Range("AC28").ClearContents()

# --- Row 29 (was row 28): "Tvåtandad spolsnäcka" record - already correct --
# (values match after the shift; nothing further to change)
